$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.714293
$ws.Range("H2").Value = 8.142879000000001
$ws.Range("I2").Value = 0.02893885961486273
$ws.Range("J2").Value = 0.02893885961486273
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 41.63852366666666
$ws.Range("N2").Value = 124.915571
$ws.Range("O2").Value = 0.08215189730289395
$ws.Range("P2").Value = 0.08215189730289395
$ws.Range("Q2").Value = 113.0191533187677
$ws.Range("R2").Value = 1017.172379868909
$ws.Range("S2").Value = 0.002377382223143068
$ws.Range("T2").Value = 0.002377382223143068

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.714293
$ws.Range("H3").Value = 8.142879000000001
$ws.Range("I3").Value = 0.02893885961486273
$ws.Range("J3").Value = 0.02893885961486273
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 211.5004576666667
$ws.Range("N3").Value = 634.5013730000001
$ws.Range("O3").Value = 0.4172857812357213
$ws.Range("P3").Value = 0.4172857812357213
$ws.Range("Q3").Value = 574.0742117414297
$ws.Range("R3").Value = 5166.667905672868
$ws.Range("S3").Value = 0.01207577464245886
$ws.Range("T3").Value = 0.01207577464245886

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.714293
$ws.Range("H4").Value = 8.142879000000001
$ws.Range("I4").Value = 0.02893885961486273
$ws.Range("J4").Value = 0.02893885961486273
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 103.0904793333333
$ws.Range("N4").Value = 309.271438
$ws.Range("O4").Value = 0.2033952629756168
$ws.Range("P4").Value = 0.2033952629756168
$ws.Range("Q4").Value = 279.8177664211113
$ws.Range("R4").Value = 2518.359897790002
$ws.Range("S4").Value = 0.005886026961579462
$ws.Range("T4").Value = 0.005886026961579462

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.714293
$ws.Range("H5").Value = 8.142879000000001
$ws.Range("I5").Value = 0.02893885961486273
$ws.Range("J5").Value = 0.02893885961486273
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 150.618525
$ws.Range("N5").Value = 451.855575
$ws.Range("O5").Value = 0.297167058485768
$ws.Range("P5").Value = 0.2971670584857679
$ws.Range("Q5").Value = 408.822808077825
$ws.Range("R5").Value = 3679.405272700425
$ws.Range("S5").Value = 0.008599675787681341
$ws.Range("T5").Value = 0.00859967578768134

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 64.641609
$ws.Range("H6").Value = 193.924827
$ws.Range("I6").Value = 0.6891866309679342
$ws.Range("J6").Value = 0.6891866309679342
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 41.63852366666666
$ws.Range("N6").Value = 124.915571
$ws.Range("O6").Value = 0.08215189730289395
$ws.Range("P6").Value = 0.08215189730289395
$ws.Range("Q6").Value = 2691.581166197913
$ws.Range("R6").Value = 24224.23049578122
$ws.Range("S6").Value = 0.0566179893298052
$ws.Range("T6").Value = 0.0566179893298052

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 64.641609
$ws.Range("H7").Value = 193.924827
$ws.Range("I7").Value = 0.6891866309679342
$ws.Range("J7").Value = 0.6891866309679342
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 211.5004576666667
$ws.Range("N7").Value = 634.5013730000001
$ws.Range("O7").Value = 0.4172857812357213
$ws.Range("P7").Value = 0.4172857812357213
$ws.Range("Q7").Value = 13671.72988780972
$ws.Range("R7").Value = 123045.5689902875
$ws.Range("S7").Value = 0.2875877817206692
$ws.Range("T7").Value = 0.2875877817206692

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 64.641609
$ws.Range("H8").Value = 193.924827
$ws.Range("I8").Value = 0.6891866309679342
$ws.Range("J8").Value = 0.6891866309679342
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 103.0904793333333
$ws.Range("N8").Value = 309.271438
$ws.Range("O8").Value = 0.2033952629756168
$ws.Range("P8").Value = 0.2033952629756168
$ws.Range("Q8").Value = 6663.934456687914
$ws.Range("R8").Value = 59975.41011019122
$ws.Range("S8").Value = 0.1401772960450024
$ws.Range("T8").Value = 0.1401772960450023

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 64.641609
$ws.Range("H9").Value = 193.924827
$ws.Range("I9").Value = 0.6891866309679342
$ws.Range("J9").Value = 0.6891866309679342
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 150.618525
$ws.Range("N9").Value = 451.855575
$ws.Range("O9").Value = 0.297167058485768
$ws.Range("P9").Value = 0.297167058485768
$ws.Range("Q9").Value = 9736.223801206726
$ws.Range("R9").Value = 87626.01421086052
$ws.Range("S9").Value = 0.2048035638724575
$ws.Range("T9").Value = 0.2048035638724575

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.357276
$ws.Range("H10").Value = 7.071828
$ws.Range("I10").Value = 0.02513246699508312
$ws.Range("J10").Value = 0.02513246699508312
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 41.63852366666666
$ws.Range("N10").Value = 124.915571
$ws.Range("O10").Value = 0.08215189730289395
$ws.Range("P10").Value = 0.08215189730289395
$ws.Range("Q10").Value = 98.15349251486533
$ws.Range("R10").Value = 883.381432633788
$ws.Range("S10").Value = 0.002064679847548441
$ws.Range("T10").Value = 0.002064679847548441

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.357276
$ws.Range("H11").Value = 7.071828
$ws.Range("I11").Value = 0.02513246699508312
$ws.Range("J11").Value = 0.02513246699508312
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 211.5004576666667
$ws.Range("N11").Value = 634.5013730000001
$ws.Range("O11").Value = 0.4172857812357213
$ws.Range("P11").Value = 0.4172857812357213
$ws.Range("Q11").Value = 498.5649528466494
$ws.Range("R11").Value = 4487.084575619844
$ws.Range("S11").Value = 0.01048742112442424
$ws.Range("T11").Value = 0.01048742112442424

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.357276
$ws.Range("H12").Value = 7.071828
$ws.Range("I12").Value = 0.02513246699508312
$ws.Range("J12").Value = 0.02513246699508312
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 103.0904793333333
$ws.Range("N12").Value = 309.271438
$ws.Range("O12").Value = 0.2033952629756168
$ws.Range("P12").Value = 0.2033952629756168
$ws.Range("Q12").Value = 243.0127127609627
$ws.Range("R12").Value = 2187.114414848664
$ws.Range("S12").Value = 0.005111824733690942
$ws.Range("T12").Value = 0.005111824733690941

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.357276
$ws.Range("H13").Value = 7.071828
$ws.Range("I13").Value = 0.02513246699508312
$ws.Range("J13").Value = 0.02513246699508312
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 150.618525
$ws.Range("N13").Value = 451.855575
$ws.Range("O13").Value = 0.297167058485768
$ws.Range("P13").Value = 0.2971670584857679
$ws.Range("Q13").Value = 355.0494341379001
$ws.Range("R13").Value = 3195.4449072411
$ws.Range("S13").Value = 0.007468541289419499
$ws.Range("T13").Value = 0.007468541289419498

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 24.08087733333333
$ws.Range("H14").Value = 72.242632
$ws.Range("I14").Value = 0.25674204242212
$ws.Range("J14").Value = 0.25674204242212
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 41.63852366666666
$ws.Range("N14").Value = 124.915571
$ws.Range("O14").Value = 0.08215189730289395
$ws.Range("P14").Value = 0.08215189730289395
$ws.Range("Q14").Value = 1002.692180758097
$ws.Range("R14").Value = 9024.229626822873
$ws.Range("S14").Value = 0.02109184590239724
$ws.Range("T14").Value = 0.02109184590239724

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 24.08087733333333
$ws.Range("H15").Value = 72.242632
$ws.Range("I15").Value = 0.25674204242212
$ws.Range("J15").Value = 0.25674204242212
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 211.5004576666667
$ws.Range("N15").Value = 634.5013730000001
$ws.Range("O15").Value = 0.4172857812357213
$ws.Range("P15").Value = 0.4172857812357213
$ws.Range("Q15").Value = 5093.11657701486
$ws.Range("R15").Value = 45838.04919313374
$ws.Range("S15").Value = 0.107134803748169
$ws.Range("T15").Value = 0.107134803748169

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 24.08087733333333
$ws.Range("H16").Value = 72.242632
$ws.Range("I16").Value = 0.25674204242212
$ws.Range("J16").Value = 0.25674204242212
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 103.0904793333333
$ws.Range("N16").Value = 309.271438
$ws.Range("O16").Value = 0.2033952629756168
$ws.Range("P16").Value = 0.2033952629756168
$ws.Range("Q16").Value = 2482.509187060535
$ws.Range("R16").Value = 22342.58268354482
$ws.Range("S16").Value = 0.05222011523534406
$ws.Range("T16").Value = 0.05222011523534405

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 24.08087733333333
$ws.Range("H17").Value = 72.242632
$ws.Range("I17").Value = 0.25674204242212
$ws.Range("J17").Value = 0.25674204242212
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 150.618525
$ws.Range("N17").Value = 451.855575
$ws.Range("O17").Value = 0.297167058485768
$ws.Range("P17").Value = 0.2971670584857679
$ws.Range("Q17").Value = 3627.0262246526
$ws.Range("R17").Value = 32643.2360218734
$ws.Range("S17").Value = 0.07629527753620964
$ws.Range("T17").Value = 0.07629527753620963
